# Scope and Limitations 1.2
# Rename the "Customers able to view all issues" row label (cell B4 on Sheet1)
# to "All concurrent users see updated issue information".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = "All concurrent users see updated issue information"
